# Release blancoRestGeneratorTs 0.1.0 !
#
# Adds a new "import文の自動生成" (import-statement auto-generation) row to
# the valueObject sheet's "common" settings block (between the existing
# "デフォルト値の変形" row and the blank spacer row before the "継承" block),
# and updates the selection/validation ranges that depend on it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("valueObject")

# --- Insert a new row above row 14 -----------------------------------
# This naturally shifts every row from 14 downwards by one (rows 14-43
# become 15-44), and shifts the merged cells / E59 data validation /
# row-relative formatting that live below the insertion point.
$ws.Rows.Item(14).Insert()

# --- Populate the new row 14 ------------------------------------------
$ws.Range("A14").Value = "import文の自動生成"
$ws.Range("B14").Value = ""
$ws.Range("C14").Value = "○"
$ws.Range("D14").Value = "/* TypeScript 独自。blancoで一括生成されたクラスについて、import文を自動生成します。 */"

# Formatting for A14 / B14, matching the look of the other rows in this
# "common" block (fontId=1 / fillId=2 / thin borders) but left-aligned.
$ws.Range("A14").Font.Name = "ＭＳ ゴシック"
$ws.Range("A14").Interior.PatternColorIndex = 27
$ws.Range("A14").Interior.ColorIndex = 42
$ws.Range("A14").Borders.Item(7).LineStyle = 1
$ws.Range("A14").Borders.Item(7).Weight = 2
$ws.Range("A14").Borders.Item(8).LineStyle = 1
$ws.Range("A14").Borders.Item(8).Weight = 2
$ws.Range("A14").Borders.Item(9).LineStyle = 1
$ws.Range("A14").Borders.Item(9).Weight = 2
$ws.Range("A14").HorizontalAlignment = -4131
$ws.Range("A14").VerticalAlignment = -4108

$ws.Range("B14").Font.Name = "ＭＳ ゴシック"
$ws.Range("B14").Interior.PatternColorIndex = 27
$ws.Range("B14").Interior.ColorIndex = 42
$ws.Range("B14").Borders.Item(9).LineStyle = 1
$ws.Range("B14").Borders.Item(9).Weight = 2
$ws.Range("B14").Borders.Item(8).LineStyle = 1
$ws.Range("B14").Borders.Item(8).Weight = 2
$ws.Range("B14").HorizontalAlignment = -4131
$ws.Range("B14").VerticalAlignment = -4108

$ws.Range("C14").Font.Name = "ＭＳ ゴシック"
$ws.Range("C14").Interior.PatternColorIndex = 9
$ws.Range("C14").Interior.ColorIndex = 26
$ws.Range("C14").Borders.Item(7).LineStyle = 1
$ws.Range("C14").Borders.Item(8).LineStyle = 1
$ws.Range("C14").Borders.Item(9).LineStyle = 1
$ws.Range("C14").Borders.Item(10).LineStyle = 1
$ws.Range("C14").HorizontalAlignment = -4108
$ws.Range("C14").VerticalAlignment = -4108

# --- Extend the "adjustDefaultValue" validation list to include the
#     newly inserted row (C12:C13 -> C12:C14). ------------------------
$ws.Range("C12:C14").Validation.Delete()
$ws.Range("C12:C14").Validation.Add(3, 1, 1, "adjustDefaultValue")

# --- Update the current selection to match the edited cell ------------
$ws.Range("D14").Select()
